$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.483.09'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +3.23%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.823.03'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +4.58%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '343.95'

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.05%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3819'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +1.27%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3529'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +4.11%  '

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.78%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.234'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +3.51%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07725'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +3.37%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.005'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.18%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.11'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +8.06%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.603'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.16%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.826.58'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +4.51%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.209'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.41%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001123'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +3.22%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06742'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.80%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '86.72'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +3.80%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.05%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.52'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +4.48%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.532'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +5.34%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '13.15'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.52%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '27.510.44'

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.480'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.24%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.684'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +8.24%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '21.95'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +11.57%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.475'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +3.49%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '153.16'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.58%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.035.98'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +4.96%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '135.28'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +2.52%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.329'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +3.26%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.090'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.84%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '13.84'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +6.29%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.08792'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +1.37%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.696'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.96%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.615'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.20%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.6985'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +11.60%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '9.146'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +6.19%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06509'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +3.09%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2253'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +3.17%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.02399'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.65%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.320'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +7.31%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '14.80'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +3.41%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6576'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +8.38%  '

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.948'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.46%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.185'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +5.53%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '133.16'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.20%  '

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.70%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '80.94'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +3.90%  '
